# Regenerate the "K" column (G) values for save_data: recompute s_vals and
# write the newly calculated K values over the old Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2 = 4;  3 = 0;  4 = 1;  5 = 1;  6 = 1;  7 = 1;  8 = 2;  9 = 1;  10 = 2;
    11 = 1; 12 = 1; 13 = 0; 14 = 2; 15 = 1; 16 = 2; 17 = 2; 18 = 2; 19 = 2; 20 = 2;
    21 = 0; 22 = 0; 23 = 1; 24 = 2; 25 = 2; 26 = 1; 27 = 1; 28 = 3; 29 = 1; 30 = 0;
    31 = 3; 32 = 2; 33 = 0; 34 = 0; 35 = 2; 36 = 1; 37 = 1; 38 = 1; 39 = 1; 40 = 2;
    41 = 1; 42 = 0; 43 = 1; 44 = 1; 45 = 1; 46 = 3; 47 = 2; 48 = 0; 49 = 1; 50 = 0;
    51 = 1; 52 = 0; 53 = 1; 54 = 0; 55 = 2; 56 = 0; 57 = 2; 58 = 1; 59 = 1; 60 = 1;
    61 = 2; 62 = 2; 63 = 1; 64 = 1; 65 = 2; 66 = 2; 67 = 1; 68 = 1; 69 = 1; 70 = 1;
    71 = 1; 72 = 3; 73 = 1; 74 = 2; 75 = 2; 76 = 3; 77 = 1; 78 = 0; 79 = 1; 80 = 1;
    81 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
